$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 128, shifting existing rows down
$ws.Rows.Item(128).Insert()

# Populate new row 128 with data (copy of original row 128's static fields + new values)
$ws.Cells.Item(128, 1).Value = 4
$ws.Cells.Item(128, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(128, 3).Value = "Los Lagos"
$ws.Cells.Item(128, 4).Value = 44704
$ws.Cells.Item(128, 5).Value = 10
$ws.Cells.Item(128, 6).Value = 100112009
$ws.Cells.Item(128, 7).Value = "Acelga"
$ws.Cells.Item(128, 8).Value = "Sin especificar"
$ws.Cells.Item(128, 9).Value = "Primera"
$ws.Cells.Item(128, 10).Value = 35
$ws.Cells.Item(128, 11).Value = 12000
$ws.Cells.Item(128, 12).Value = 12000
$ws.Cells.Item(128, 13).Value = 12000
$ws.Cells.Item(128, 14).Value = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(128, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(128, 16).Value = 1000
$ws.Cells.Item(128, 17).Value = 12
$ws.Cells.Item(128, 18).Value = "Hortaliza"

Write-Host "Inserted new row 128; sheet now has" $ws.UsedRange.Rows.Count "rows"
